$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 351.875
$ws.Range("I12").Value = 236.5
$ws.Range("K12").Value = 236.5
$ws.Range("M12").Value = -66.5

$ws.Range("H40").Value = 8144.4443
$ws.Range("I40").Value = 8715.666999999999
$ws.Range("K40").Value = 8715.666999999999
$ws.Range("M40").Value = -8540.666999999999

$ws.Range("H70").Value = 3267.4285
$ws.Range("J70").Value = 3561.75
$ws.Range("L70").Value = 10685.25
$ws.Range("N70").Value = -11225.25

$ws.Range("H73").Value = 3267.4285
$ws.Range("J73").Value = 3561.75
$ws.Range("L73").Value = 10685.25
$ws.Range("N73").Value = -12557.25

$ws.Range("H107").Value = 405.75
$ws.Range("I107").Value = 405.75
$ws.Range("K107").Value = 405.75
$ws.Range("M107").Value = 1514.25

$ws.Range("H127").Value = 891.6667
$ws.Range("I127").Value = 615.4286
$ws.Range("J127").Value = 1858.5
$ws.Range("K127").Value = 1846.2858
$ws.Range("L127").Value = 5575.5
$ws.Range("M127").Value = 3113.7142
$ws.Range("N127").Value = -15495.5

$ws.Range("H132").Value = 6681.636
$ws.Range("I132").Value = 6849.8
$ws.Range("K132").Value = 20549.4
$ws.Range("M132").Value = -18019.4

$ws.Range("H138").Value = 6471.3096
$ws.Range("I138").Value = 4195.6
$ws.Range("J138").Value = 6778.838
$ws.Range("K138").Value = 12586.8
$ws.Range("L138").Value = 20336.514
$ws.Range("M138").Value = -7446.800000000001
$ws.Range("N138").Value = -30616.514

$ws.Range("H141").Value = 1580
$ws.Range("I141").Value = 1514.2858
$ws.Range("K141").Value = 4542.857400000001
$ws.Range("M141").Value = 637.1425999999992

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2321.2979
$ws.Range("I32").Value = 2401.721
$ws.Range("J32").Value = 1456.75
$ws.Range("K32").Value = 2401.721
$ws.Range("L32").Value = 1456.75
$ws.Range("M32").Value = -2114.721
$ws.Range("N32").Value = -2030.75

$ws.Range("H33").Value = 4666
$ws.Range("I33").Value = 4000
$ws.Range("J33").Value = 4999
$ws.Range("K33").Value = 4000
$ws.Range("L33").Value = 4999
$ws.Range("M33").Value = -3671
$ws.Range("N33").Value = -5657

$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("M36").ClearContents()

$ws.Range("H61").Value = 2059.7222
$ws.Range("I61").Value = 1648.1154
$ws.Range("J61").Value = 3129.9
$ws.Range("K61").Value = 1648.1154
$ws.Range("L61").Value = 3129.9
$ws.Range("M61").Value = -1436.1154
$ws.Range("N61").Value = -3553.9

$ws.Range("H132").Value = 2851203.2
$ws.Range("I132").Value = 4275248.5
$ws.Range("J132").Value = 3112.6667
$ws.Range("K132").Value = 12825745.5
$ws.Range("L132").Value = 9338.000100000001
$ws.Range("M132").Value = -12823215.5
$ws.Range("N132").Value = -14398.0001

$ws.Range("H136").Value = 2059.7222
$ws.Range("I136").Value = 1648.1154
$ws.Range("J136").Value = 3129.9
$ws.Range("K136").Value = 4944.3462
$ws.Range("L136").Value = 9389.700000000001
$ws.Range("M136").Value = -2394.3462
$ws.Range("N136").Value = -14489.7

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 5000
$ws.Range("I7").Value = 3333.3333
$ws.Range("K7").Value = 3333.3333
$ws.Range("M7").Value = -3220.3333

$ws.Range("H105").Value = 32502736
$ws.Range("I105").Value = 3334965.2
$ws.Range("K105").Value = 3334965.2
$ws.Range("M105").Value = -3333218.2

$ws.Range("H134").Value = 2334.6086
$ws.Range("I134").Value = 1552.4
$ws.Range("J134").Value = 2936.3076
$ws.Range("K134").Value = 4657.200000000001
$ws.Range("L134").Value = 8808.9228
$ws.Range("M134").Value = -2122.200000000001
$ws.Range("N134").Value = -13878.9228

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 6099.8423
$ws.Range("I58").Value = 3928
$ws.Range("J58").Value = 7366.75
$ws.Range("K58").Value = 3928
$ws.Range("L58").Value = 7366.75
$ws.Range("M58").Value = -3725
$ws.Range("N58").Value = -7772.75

$ws.Range("H134").Value = 3519.5518
$ws.Range("I134").Value = 3367.0908
$ws.Range("K134").Value = 10101.2724
$ws.Range("M134").Value = -7566.2724

$ws.Range("H136").Value = 6099.8423
$ws.Range("I136").Value = 3928
$ws.Range("J136").Value = 7366.75
$ws.Range("K136").Value = 11784
$ws.Range("L136").Value = 22100.25
$ws.Range("M136").Value = -9234
$ws.Range("N136").Value = -27200.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 3719.4443
$ws.Range("I8").Value = 3719.4443
$ws.Range("K8").Value = 11158.3329
$ws.Range("M8").Value = -11019.3329

$ws.Range("H23").Value = 4363.391
$ws.Range("I23").Value = 593.1429000000001
$ws.Range("J23").Value = 6012.875
$ws.Range("K23").Value = 1779.4287
$ws.Range("L23").Value = 18038.625
$ws.Range("M23").Value = -1544.4287
$ws.Range("N23").Value = -18508.625

$ws.Range("H64").Value = 1850
$ws.Range("J64").Value = 2000
$ws.Range("L64").Value = 6000
$ws.Range("N64").Value = -6540

$ws.Range("H67").Value = 1850
$ws.Range("J67").Value = 2000
$ws.Range("L67").Value = 6000
$ws.Range("N67").Value = -7872

$ws.Range("H68").Value = 4447939.5
$ws.Range("J68").Value = 4003874.2
$ws.Range("L68").Value = 12011622.6
$ws.Range("N68").Value = -12013244.6

$ws.Range("H71").Value = 4447939.5
$ws.Range("J71").Value = 4003874.2
$ws.Range("L71").Value = 36034867.8
$ws.Range("N71").Value = -36042979.8

$ws.Range("H98").Value = 421.6154
$ws.Range("I98").Value = 408.33334
$ws.Range("J98").Value = 425.6
$ws.Range("K98").Value = 1225.00002
$ws.Range("L98").Value = 1276.8
$ws.Range("M98").Value = 272.9999800000001
$ws.Range("N98").Value = -4272.8

$ws.Range("H114").Value = 1821.7646
$ws.Range("J114").Value = 2552.3333
$ws.Range("L114").Value = 7656.999899999999
$ws.Range("N114").Value = -14164.9999

$ws.Range("H117").Value = 2878.5
$ws.Range("J117").Value = 2881
$ws.Range("L117").Value = 8643
$ws.Range("N117").Value = -15527

$ws.Range("H121").Value = 5614658.5
$ws.Range("I121").Value = 14286458
$ws.Range("J121").Value = 96240.82000000001
$ws.Range("K121").Value = 42859374
$ws.Range("L121").Value = 288722.46
$ws.Range("M121").Value = -42858064
$ws.Range("N121").Value = -291342.46

$ws.Range("H131").Value = 7791.909
$ws.Range("J131").Value = 2117.5557
$ws.Range("L131").Value = 6352.6671
$ws.Range("N131").Value = -16432.6671

$ws.Range("H139").Value = 2459.5715
$ws.Range("I139").Value = 689.3570999999999
$ws.Range("K139").Value = 2068.0713
$ws.Range("M139").Value = 3071.9287

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 651.37036
$ws.Range("I97").Value = 600.3333
$ws.Range("K97").Value = 600.3333
$ws.Range("M97").Value = -104.3333

$ws.Range("H122").Value = 76927120
$ws.Range("I122").Value = 76924690
$ws.Range("J122").Value = 76929544
$ws.Range("K122").Value = 230774070
$ws.Range("L122").Value = 230788632
$ws.Range("M122").Value = -230771620
$ws.Range("N122").Value = -230793532

$ws.Range("H132").Value = 2247.641
$ws.Range("I132").Value = 2399.5417
$ws.Range("J132").Value = 2004.6
$ws.Range("K132").Value = 7198.625100000001
$ws.Range("L132").Value = 6013.799999999999
$ws.Range("M132").Value = -4668.625100000001
$ws.Range("N132").Value = -11073.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2405.3333
$ws.Range("I7").Value = 2186.4
$ws.Range("K7").Value = 2186.4
$ws.Range("M7").Value = -2074.4

$ws.Range("H16").Value = 3495
$ws.Range("I16").Value = 3495
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 3495
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -3325
$ws.Range("N16").ClearContents()

$ws.Range("H46").Value = 1666.3334
$ws.Range("I46").Value = 999.5
$ws.Range("K46").Value = 999.5
$ws.Range("M46").Value = -811.5

$ws.Range("H56").Value = 15209.6
$ws.Range("I56").Value = 525
$ws.Range("J56").Value = 24999.334
$ws.Range("K56").Value = 525
$ws.Range("L56").Value = 24999.334
$ws.Range("M56").Value = 166
$ws.Range("N56").Value = -26381.334

$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").ClearContents()

$ws.Range("H126").Value = 2405.3333
$ws.Range("I126").Value = 2186.4
$ws.Range("K126").Value = 6559.200000000001
$ws.Range("M126").Value = -4089.200000000001

$ws.Range("H132").Value = 9699.1
$ws.Range("J132").Value = 7199.8
$ws.Range("L132").Value = 21599.4
$ws.Range("N132").Value = -26659.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H7").Value = 2577.5
$ws.Range("I7").Value = 2750
$ws.Range("J7").Value = 2405
$ws.Range("K7").Value = 2750
$ws.Range("L7").Value = 2405
$ws.Range("M7").Value = -2637
$ws.Range("N7").Value = -2631

$ws.Range("H9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("N9").ClearContents()

$ws.Range("H14").Value = 9202
$ws.Range("I14").Value = 9000
$ws.Range("K14").Value = 9000
$ws.Range("M14").Value = -8832

$ws.Range("H96").Value = 1681.9
$ws.Range("I96").Value = 1604
$ws.Range("J96").Value = 1759.8
$ws.Range("K96").Value = 1604
$ws.Range("L96").Value = 1759.8
$ws.Range("M96").Value = -231
$ws.Range("N96").Value = -4505.8

$ws.Range("H132").Value = 2666.6458
$ws.Range("I132").Value = 2614.6428
$ws.Range("K132").Value = 7843.928400000001
$ws.Range("M132").Value = -5313.928400000001

$ws.Range("H136").Value = 5382806
$ws.Range("J136").Value = 7367.6875
$ws.Range("L136").Value = 22103.0625
$ws.Range("N136").Value = -27203.0625
